$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in computed metrics for rows 2-7 (C:G), previously blank
$data = @{
    2 = @(0.626, 14.996, 12.95, 235.891, 6)
    3 = @(-1.887, 36.077, 28.133, 234.863, 6)
    4 = @(-1.331, 38.377, 28.953, 238.204, 6)
    5 = @(0.653, 14.442, 12.702, 237.873, 6)
    6 = @(-1.969, 36.583, 27.725, 234.597, 6)
    7 = @(-1.383, 38.805, 28.893, 240.171, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}

# Clear out N_obs_test (column G) values for rows 8-10, leaving them blank
$ws.Range("G8:G10").ClearContents()
